# Update iserv_stats for 2025-09 (row 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6293
$ws.Range("C22").Value = 991
$ws.Range("D22").Value = 5816328
$ws.Range("E22").Value = 924.2536151279199
$ws.Range("F22").Value = 8.331898777758639
$ws.Range("G22").Value = 3.661087866108792
$ws.Range("H22").Value = 26.4870084116946
